$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 11111
$ws.Range("I6").Value = 11111
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 33333
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -33221
$ws.Range("N6").ClearContents()
$ws.Range("H32").Value = 1867.6666
$ws.Range("I32").Value = 1750
$ws.Range("J32").Value = 1926.5
$ws.Range("K32").Value = 1750
$ws.Range("L32").Value = 1926.5
$ws.Range("M32").Value = -1424
$ws.Range("N32").Value = -2578.5
$ws.Range("H39").Value = 785.5
$ws.Range("I39").Value = 804.5333000000001
$ws.Range("K39").Value = 2413.5999
$ws.Range("M39").Value = -2117.5999
$ws.Range("H107").Value = 1271443.5
$ws.Range("I107").Value = 1779150.9
$ws.Range("J107").Value = 2175
$ws.Range("K107").Value = 1779150.9
$ws.Range("L107").Value = 2175
$ws.Range("M107").Value = -1777230.9
$ws.Range("N107").Value = -6015
$ws.Range("H133").Value = 16833.334
$ws.Range("J133").Value = 16833.334
$ws.Range("L133").Value = 16833.334
$ws.Range("N133").Value = -26953.334
$ws.Range("H137").Value = 1303.525
$ws.Range("I137").Value = 864.9231
$ws.Range("J137").Value = 1514.7037
$ws.Range("K137").Value = 2594.7693
$ws.Range("L137").Value = 4544.1111
$ws.Range("M137").Value = -44.76929999999993
$ws.Range("N137").Value = -9644.1111
$ws.Range("H138").Value = 5156922
$ws.Range("I138").Value = 1966.2
$ws.Range("J138").Value = 7465111.5
$ws.Range("K138").Value = 5898.6
$ws.Range("L138").Value = 22395334.5
$ws.Range("M138").Value = -758.6000000000004
$ws.Range("N138").Value = -22405614.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2449.1667
$ws.Range("I45").Value = 1539
$ws.Range("K45").Value = 1539
$ws.Range("M45").Value = -1162
$ws.Range("H74").Value = 1588.8148
$ws.Range("I74").Value = 1895.6
$ws.Range("J74").Value = 1205.3334
$ws.Range("K74").Value = 1895.6
$ws.Range("L74").Value = 1205.3334
$ws.Range("M74").Value = -1021.6
$ws.Range("N74").Value = -2953.3334
$ws.Range("H77").Value = 1588.8148
$ws.Range("I77").Value = 1895.6
$ws.Range("J77").Value = 1205.3334
$ws.Range("K77").Value = 9478
$ws.Range("L77").Value = 6026.666999999999
$ws.Range("M77").Value = -5110
$ws.Range("N77").Value = -14762.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1425.7333
$ws.Range("I20").Value = 1125.4706
$ws.Range("J20").Value = 1818.3846
$ws.Range("K20").Value = 1125.4706
$ws.Range("L20").Value = 1818.3846
$ws.Range("M20").Value = -878.4706000000001
$ws.Range("N20").Value = -2312.3846
$ws.Range("H59").Value = 44500
$ws.Range("J59").Value = 44500
$ws.Range("L59").Value = 44500
$ws.Range("N59").Value = -46194
$ws.Range("H64").Value = 940.375
$ws.Range("I64").Value = 756
$ws.Range("K64").Value = 756
$ws.Range("M64").Value = -531
$ws.Range("H67").Value = 940.375
$ws.Range("I67").Value = 756
$ws.Range("K67").Value = 756
$ws.Range("M67").Value = 24
$ws.Range("H133").Value = 49000
$ws.Range("J133").Value = 49000
$ws.Range("L133").Value = 49000
$ws.Range("N133").Value = -59120
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1299.9565
$ws.Range("I31").Value = 952.93475
$ws.Range("J31").Value = 1646.9783
$ws.Range("K31").Value = 952.93475
$ws.Range("L31").Value = 1646.9783
$ws.Range("M31").Value = -657.93475
$ws.Range("N31").Value = -2236.9783
$ws.Range("H34").Value = 1299.9565
$ws.Range("I34").Value = 952.93475
$ws.Range("J34").Value = 1646.9783
$ws.Range("K34").Value = 952.93475
$ws.Range("L34").Value = 1646.9783
$ws.Range("M34").Value = -750.93475
$ws.Range("N34").Value = -2050.9783
$ws.Range("H58").Value = 2117.0527
$ws.Range("I58").Value = 854.8889
$ws.Range("K58").Value = 854.8889
$ws.Range("M58").Value = -651.8889
$ws.Range("H107").Value = 552.2917
$ws.Range("I107").Value = 446
$ws.Range("K107").Value = 446
$ws.Range("M107").Value = 1474
$ws.Range("H132").Value = 2658.875
$ws.Range("I132").Value = 1933.6818
$ws.Range("J132").Value = 4254.3
$ws.Range("K132").Value = 5801.0454
$ws.Range("L132").Value = 12762.9
$ws.Range("M132").Value = -3271.0454
$ws.Range("N132").Value = -17822.9
$ws.Range("H136").Value = 2117.0527
$ws.Range("I136").Value = 854.8889
$ws.Range("K136").Value = 2564.6667
$ws.Range("M136").Value = -14.66670000000022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1359.2858
$ws.Range("I5").Value = 614
$ws.Range("J5").Value = 2477.2144
$ws.Range("K5").Value = 1842
$ws.Range("L5").Value = 7431.6432
$ws.Range("M5").Value = -1730
$ws.Range("N5").Value = -7655.6432
$ws.Range("H7").Value = 124.25
$ws.Range("I7").Value = 115.666664
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 346.999992
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = -234.999992
$ws.Range("N7").Value = -674
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H80").Value = 1200
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1200
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3600
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5472
$ws.Range("H83").Value = 1200
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1200
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 10800
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -20160
$ws.Range("H92").Value = 740
$ws.Range("J92").Value = 763.3333
$ws.Range("L92").Value = 2289.9999
$ws.Range("N92").Value = -4785.9999
$ws.Range("H102").Value = 4029
$ws.Range("J102").Value = 4029
$ws.Range("L102").Value = 12087
$ws.Range("N102").Value = -16955
$ws.Range("H107").Value = 995.4032
$ws.Range("I107").Value = 324.73685
$ws.Range("J107").Value = 1291.7441
$ws.Range("K107").Value = 974.21055
$ws.Range("L107").Value = 3875.2323
$ws.Range("M107").Value = 945.78945
$ws.Range("N107").Value = -7715.2323
$ws.Range("H132").Value = 1172
$ws.Range("I132").Value = 867.6667
$ws.Range("J132").Value = 1232.8667
$ws.Range("K132").Value = 7809.0003
$ws.Range("L132").Value = 11095.8003
$ws.Range("M132").Value = -5279.0003
$ws.Range("N132").Value = -16155.8003
$ws.Range("H135").Value = 1359.2858
$ws.Range("I135").Value = 614
$ws.Range("J135").Value = 2477.2144
$ws.Range("K135").Value = 5526
$ws.Range("L135").Value = 22294.9296
$ws.Range("M135").Value = -2991
$ws.Range("N135").Value = -27364.9296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 64466.668
$ws.Range("J138").Value = 64466.668
$ws.Range("L138").Value = 64466.668
$ws.Range("N138").Value = -74746.66800000001
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3504.5925
$ws.Range("I40").Value = 1977.4
$ws.Range("J40").Value = 4402.9414
$ws.Range("K40").Value = 1977.4
$ws.Range("L40").Value = 4402.9414
$ws.Range("M40").Value = -1841.4
$ws.Range("N40").Value = -4674.9414
$ws.Range("H55").Value = 402
$ws.Range("I55").Value = 383.6
$ws.Range("J55").Value = 420.4
$ws.Range("K55").Value = 383.6
$ws.Range("L55").Value = 420.4
$ws.Range("M55").Value = -210.6
$ws.Range("N55").Value = -766.4
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 70000
$ws.Range("I141").Value = 70000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 70000
$ws.Range("L141").ClearContents()
$ws.Range("M141").Value = -64820
$ws.Range("N141").ClearContents()
